$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "'301.32"
$ws.Range("E2").Value = "'0.03%"
$ws.Range("D3").Value = "'32.28"
$ws.Range("E3").Value = "'1.84%"
$ws.Range("D4").Value = "'5.019"
$ws.Range("E4").Value = "'-1.47%"
$ws.Range("D5").Value = "'0.07629"
$ws.Range("D6").Value = "'1.952"
$ws.Range("E6").Value = "'-12.77%"
$ws.Range("D7").Value = "'7.872"
$ws.Range("E7").Value = "'1.01%"
$ws.Range("D8").Value = "'3.782"
$ws.Range("E8").Value = "'-0.95%"
$ws.Range("D9").Value = "'0.9181"
$ws.Range("E9").Value = "'-0.04%"
$ws.Range("D10").Value = "'0.1758"
$ws.Range("E10").Value = "'-0.03%"
$ws.Range("D11").Value = "'0.07829"
$ws.Range("E11").Value = "'3.77%"
$ws.Range("D12").Value = "'0.08510"
$ws.Range("E12").Value = "'-5.31%"
$ws.Range("D13").Value = "'0.03161"
$ws.Range("E13").Value = "'4.35%"
$ws.Range("D14").Value = "'0.1000"
$ws.Range("E14").Value = "'-0.34%"
$ws.Range("D15").Value = "'0.001509"
$ws.Range("E15").Value = "'0.25%"
$ws.Range("D16").Value = "'0.005790"
$ws.Range("E16").Value = "'-3.58%"
$ws.Range("E18").Value = "'-0.20%"
$ws.Range("E19").Value = "'-4.41%"
$ws.Range("D20").Value = "'0.3343"
$ws.Range("E20").Value = "'1.54%"
$ws.Range("E21").Value = "'-2.80%"
$ws.Range("D22").Value = "'4.265"
$ws.Range("E22").Value = "'0.63%"
$ws.Range("E23").Value = "'9.60%"
$ws.Range("D24").Value = "'0.04493"
$ws.Range("E24").Value = "'-2.07%"
$ws.Range("E25").Value = "'-2.33%"
$ws.Range("D26").Value = "'0.004393"
$ws.Range("E26").Value = "'-1.78%"
$ws.Range("E27").Value = "'0.06%"
$ws.Range("D39").Value = "'0.01700"
$ws.Range("E39").Value = "'-4.10%"
$ws.Range("D40").Value = "'0.04671"
$ws.Range("E40").Value = "'-2.35%"
$ws.Range("D41").Value = "'0.007453"
$ws.Range("E41").Value = "'0.57%"
$ws.Range("E42").Value = "'-0.81%"
$ws.Range("D43").Value = "'0.002330"
$ws.Range("E43").Value = "'6.45%"
$ws.Range("E44").Value = "'2.19%"
$ws.Range("D45").Value = "'0.00006245"
$ws.Range("E45").Value = "'0.21%"
$ws.Range("E46").Value = "'0.12%"
$ws.Range("D48").Value = "'0.8232"
$ws.Range("E48").Value = "'12.63%"
$ws.Range("E49").Value = "'0.12%"
$ws.Range("E50").Value = "'0.12%"
